$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Bump the date in A1 by one day (45310 -> 45311)
$ws.Range("A1").Value = [DateTime]::FromOADate(45311)

# Update price values
$ws.Range("D35").Value = 15097
$ws.Range("D36").Value = 3375
